# Applies the "Updated symbol list" edit (Mon Dec 12 16:55:16 UTC 2022):
# refreshed coin prices in column D, and swapped the FTXToken/GateToken rows
# (6 and 7) back to their correct Coin/Link/Volume-label values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that are stored as TEXT in the
# workbook (t="inlineStr"). Setting a plain numeric-looking string via
# .Value would be auto-coerced to a Number by the input parser, so we
# momentarily force a Text number format, assign the literal string, then
# clear the format back to the original (General / no explicit style) so
# no stray style is left behind on the cell.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" '277.25'
Set-TextValue "D3" '20.99'
Set-TextValue "D4" '6.222'
Set-TextValue "D5" '0.06193'
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue "D6" '1.571'
$ws.Range("E6").Value = '5FTXTokenFTT'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D7" '3.578'
$ws.Range("E7").Value = '6GateTokenGT'
Set-TextValue "D8" '6.561'
Set-TextValue "D9" '0.8220'
Set-TextValue "D11" '0.08216'
Set-TextValue "D12" '0.03482'
Set-TextValue "D13" '0.03101'
Set-TextValue "D15" '3.764'
Set-TextValue "D16" '0.001617'
Set-TextValue "D17" '0.04683'
Set-TextValue "D18" '0.006386'
Set-TextValue "D19" '0.006142'
Set-TextValue "D22" '3.769'
Set-TextValue "D24" '0.01388'
Set-TextValue "D25" '0.3283'
Set-TextValue "D28" '0.0002737'
Set-TextValue "D40" '0.04685'
Set-TextValue "D41" '0.007022'
Set-TextValue "D42" '0.004701'
Set-TextValue "D43" '0.1103'
Set-TextValue "D44" '0.01082'
Set-TextValue "D45" '0.00006393'
Set-TextValue "D47" '0.8454'
Set-TextValue "D48" '0.001387'
